$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures.
# D-column cells whose new value is a valid numeric literal must be
# force-typed as text (NumberFormat "@" + Style "Normal" reset) so they
# keep behaving as plain text strings, matching the source data feed.

$ws.Range("D2").Value = "67.656.11"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.478.51"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  +5.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.475.38"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.60%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "4.077.91"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "67.632.37"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000177"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "3.474.98"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.37%  "
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0717"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "2.744.91"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "329.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("E51").Value = "  -2.09%  "
